$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The prime-attribute test case was using an invalid attribute array;
# switch it to a valid selection based on the standard array (3, 2, 1, 0, -1, -2, ...).
$ws.Range("F3").Value = 1
$ws.Range("H3").Value = -2

# Leave the sheet's selection where the author left it after editing.
$ws.Range("G4").Select()
